$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.059.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.853.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5069"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3885"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08214"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.186"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.851.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.164"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001091"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06637"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.899"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.094.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.237"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.059.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.399"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1052"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.028"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.812"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.592"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02413"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06420"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.036"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2157"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.246"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6392"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.172"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.911"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5980"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.273"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.650"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.993"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.197"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06852"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
